$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 362.5
$ws.Range("I18").Value = 362.5
$ws.Range("K18").Value = 362.5
$ws.Range("M18").Value = -78.5
$ws.Range("H98").Value = 803
$ws.Range("I98").Value = 803
$ws.Range("K98").Value = 803
$ws.Range("M98").Value = 695
$ws.Range("H106").Value = 2593.9092
$ws.Range("I106").Value = 2593.9092
$ws.Range("K106").Value = 2593.9092
$ws.Range("M106").Value = -1962.9092
$ws.Range("H122").Value = 803
$ws.Range("I122").Value = 803
$ws.Range("K122").Value = 2409
$ws.Range("M122").Value = 41
$ws.Range("H125").Value = 1511.75
$ws.Range("I125").Value = 1491.5
$ws.Range("K125").Value = 13423.5
$ws.Range("M125").Value = -10963.5
$ws.Range("H137").Value = 3139.889
$ws.Range("J137").Value = 3424.4443
$ws.Range("L137").Value = 10273.3329
$ws.Range("N137").Value = -15373.3329
$ws.Range("H138").Value = 2705.2307
$ws.Range("I138").Value = 1328.9166
$ws.Range("K138").Value = 3986.7498
$ws.Range("M138").Value = 1153.2502
$ws.Range("H141").Value = 4622
$ws.Range("I141").Value = 4595
$ws.Range("K141").Value = 13785
$ws.Range("M141").Value = -8605

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9846.516
$ws.Range("I2").Value = 14768.952
$ws.Range("J2").Value = 1232.25
$ws.Range("K2").Value = 14768.952
$ws.Range("L2").Value = 1232.25
$ws.Range("M2").Value = -14655.952
$ws.Range("N2").Value = -1458.25
$ws.Range("H45").Value = 4949.9375
$ws.Range("I45").Value = 3791.875
$ws.Range("J45").Value = 6108
$ws.Range("K45").Value = 3791.875
$ws.Range("L45").Value = 6108
$ws.Range("M45").Value = -3414.875
$ws.Range("N45").Value = -6862
$ws.Range("H61").Value = 3238.611
$ws.Range("I61").Value = 1020
$ws.Range("J61").Value = 4091.923
$ws.Range("K61").Value = 1020
$ws.Range("L61").Value = 4091.923
$ws.Range("M61").Value = -808
$ws.Range("N61").Value = -4515.923
$ws.Range("H74").Value = 1987.2222
$ws.Range("I74").Value = 2029.4375
$ws.Range("J74").Value = 1649.5
$ws.Range("K74").Value = 2029.4375
$ws.Range("L74").Value = 1649.5
$ws.Range("M74").Value = -1155.4375
$ws.Range("N74").Value = -3397.5
$ws.Range("H77").Value = 1987.2222
$ws.Range("I77").Value = 2029.4375
$ws.Range("J77").Value = 1649.5
$ws.Range("K77").Value = 10147.1875
$ws.Range("L77").Value = 8247.5
$ws.Range("M77").Value = -5779.1875
$ws.Range("N77").Value = -16983.5
$ws.Range("H86").Value = 50314
$ws.Range("J86").Value = 50314
$ws.Range("L86").Value = 50314
$ws.Range("N86").Value = -52686
$ws.Range("H89").Value = 50314
$ws.Range("J89").Value = 50314
$ws.Range("L89").Value = 150942
$ws.Range("N89").Value = -162798
$ws.Range("H97").Value = 928.5833
$ws.Range("I97").Value = 928.5833
$ws.Range("K97").Value = 928.5833
$ws.Range("M97").Value = -432.5833
$ws.Range("H116").Value = 9846.516
$ws.Range("I116").Value = 14768.952
$ws.Range("J116").Value = 1232.25
$ws.Range("K116").Value = 14768.952
$ws.Range("L116").Value = 1232.25
$ws.Range("M116").Value = -12474.952
$ws.Range("N116").Value = -5820.25
$ws.Range("H122").Value = 4266.706
$ws.Range("I122").Value = 2041.64
$ws.Range("J122").Value = 10447.444
$ws.Range("K122").Value = 6124.92
$ws.Range("L122").Value = 31342.332
$ws.Range("M122").Value = -3674.92
$ws.Range("N122").Value = -36242.33199999999
$ws.Range("H132").Value = 3346.1135
$ws.Range("I132").Value = 3292.1292
$ws.Range("J132").Value = 3474.8462
$ws.Range("K132").Value = 9876.3876
$ws.Range("L132").Value = 10424.5386
$ws.Range("M132").Value = -7346.3876
$ws.Range("N132").Value = -15484.5386
$ws.Range("H136").Value = 3238.611
$ws.Range("I136").Value = 1020
$ws.Range("J136").Value = 4091.923
$ws.Range("K136").Value = 3060
$ws.Range("L136").Value = 12275.769
$ws.Range("M136").Value = -510
$ws.Range("N136").Value = -17375.769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9846.516
$ws.Range("I3").Value = 14768.952
$ws.Range("J3").Value = 1232.25
$ws.Range("K3").Value = 14768.952
$ws.Range("L3").Value = 1232.25
$ws.Range("M3").Value = -14654.952
$ws.Range("N3").Value = -1460.25
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H75").Value = 12299.363
$ws.Range("I75").Value = 7254.8887
$ws.Range("K75").Value = 7254.8887
$ws.Range("M75").Value = -6318.8887
$ws.Range("H78").Value = 12299.363
$ws.Range("I78").Value = 7254.8887
$ws.Range("K78").Value = 21764.6661
$ws.Range("M78").Value = -17084.6661
$ws.Range("H86").Value = 5739
$ws.Range("I86").Value = 1800
$ws.Range("J86").Value = 7427.143
$ws.Range("K86").Value = 1800
$ws.Range("L86").Value = 7427.143
$ws.Range("M86").Value = -677
$ws.Range("N86").Value = -9673.143
$ws.Range("H89").Value = 5739
$ws.Range("I89").Value = 1800
$ws.Range("J89").Value = 7427.143
$ws.Range("K89").Value = 9000
$ws.Range("L89").Value = 37135.715
$ws.Range("M89").Value = -3384
$ws.Range("N89").Value = -48367.715
$ws.Range("H94").Value = 3937.1428
$ws.Range("I94").Value = 2926.6667
$ws.Range("J94").Value = 10000
$ws.Range("K94").Value = 2926.6667
$ws.Range("L94").Value = 10000
$ws.Range("M94").Value = -2475.6667
$ws.Range("N94").Value = -10902
$ws.Range("H99").Value = 24049.611
$ws.Range("I99").Value = 28107.334
$ws.Range("J99").Value = 3761
$ws.Range("K99").Value = 28107.334
$ws.Range("L99").Value = 3761
$ws.Range("M99").Value = -26609.334
$ws.Range("N99").Value = -6757
$ws.Range("H134").Value = 2567.9092
$ws.Range("I134").Value = 1290.697
$ws.Range("K134").Value = 3872.090999999999
$ws.Range("M134").Value = -1337.090999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 239
$ws.Range("I7").Value = 17.8
$ws.Range("J7").Value = 397
$ws.Range("K7").Value = 17.8
$ws.Range("L7").Value = 397
$ws.Range("M7").Value = 95.2
$ws.Range("N7").Value = -623
$ws.Range("H16").Value = 2005.5
$ws.Range("I16").Value = 1674
$ws.Range("K16").Value = 1674
$ws.Range("M16").Value = -1387
$ws.Range("H31").Value = 2809.4878
$ws.Range("I31").Value = 1178.9166
$ws.Range("J31").Value = 5111.4707
$ws.Range("K31").Value = 1178.9166
$ws.Range("L31").Value = 5111.4707
$ws.Range("M31").Value = -883.9166
$ws.Range("N31").Value = -5701.4707
$ws.Range("H34").Value = 2809.4878
$ws.Range("I34").Value = 1178.9166
$ws.Range("J34").Value = 5111.4707
$ws.Range("K34").Value = 1178.9166
$ws.Range("L34").Value = 5111.4707
$ws.Range("M34").Value = -976.9166
$ws.Range("N34").Value = -5515.4707
$ws.Range("H58").Value = 1622.2858
$ws.Range("I58").Value = 1251.5238
$ws.Range("J58").Value = 2734.5715
$ws.Range("K58").Value = 1251.5238
$ws.Range("L58").Value = 2734.5715
$ws.Range("M58").Value = -1048.5238
$ws.Range("N58").Value = -3140.5715
$ws.Range("H105").Value = 1013
$ws.Range("J105").Value = 907.5
$ws.Range("L105").Value = 907.5
$ws.Range("N105").Value = -4401.5
$ws.Range("H113").Value = 2005.5
$ws.Range("I113").Value = 1674
$ws.Range("K113").Value = 1674
$ws.Range("M113").Value = 496
$ws.Range("H122").Value = 1279189.9
$ws.Range("I122").Value = 1703420
$ws.Range("J122").Value = 6499.5
$ws.Range("K122").Value = 5110260
$ws.Range("L122").Value = 19498.5
$ws.Range("M122").Value = -5107810
$ws.Range("N122").Value = -24398.5
$ws.Range("H132").Value = 5130.647
$ws.Range("I132").Value = 4099.5
$ws.Range("K132").Value = 12298.5
$ws.Range("M132").Value = -9768.5
$ws.Range("H136").Value = 1622.2858
$ws.Range("I136").Value = 1251.5238
$ws.Range("J136").Value = 2734.5715
$ws.Range("K136").Value = 3754.5714
$ws.Range("L136").Value = 8203.7145
$ws.Range("M136").Value = -1204.5714
$ws.Range("N136").Value = -13303.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 41674480
$ws.Range("I22").Value = 699.3333
$ws.Range("J22").Value = 66678748
$ws.Range("K22").Value = 2097.9999
$ws.Range("L22").Value = 200036244
$ws.Range("M22").Value = -1928.9999
$ws.Range("N22").Value = -200036582
$ws.Range("H27").Value = 41674480
$ws.Range("I27").Value = 699.3333
$ws.Range("J27").Value = 66678748
$ws.Range("K27").Value = 2097.9999
$ws.Range("L27").Value = 200036244
$ws.Range("M27").Value = -1995.9999
$ws.Range("N27").Value = -200036448
$ws.Range("H80").Value = 4599
$ws.Range("I80").Value = 3998.75
$ws.Range("K80").Value = 11996.25
$ws.Range("M80").Value = -11060.25
$ws.Range("H83").Value = 4599
$ws.Range("I83").Value = 3998.75
$ws.Range("K83").Value = 35988.75
$ws.Range("M83").Value = -31308.75
$ws.Range("H121").Value = 2100.5518
$ws.Range("J121").Value = 2598.739
$ws.Range("L121").Value = 7796.217000000001
$ws.Range("N121").Value = -10416.217
$ws.Range("H122").Value = 1634.08
$ws.Range("I122").Value = 1317.7
$ws.Range("K122").Value = 11859.3
$ws.Range("M122").Value = -9409.300000000001
$ws.Range("H128").Value = 999999.7
$ws.Range("I128").Value = 999999.7
$ws.Range("K128").Value = 2999999.1
$ws.Range("M128").Value = -2995019.1
$ws.Range("H136").Value = 2424.9524
$ws.Range("I136").Value = 966.63635
$ws.Range("K136").Value = 2899.90905
$ws.Range("M136").Value = 2200.09095

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 10000
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H50").Value = 10000
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H70").Value = 107192.37
$ws.Range("I70").Value = 189534.33
$ws.Range("K70").Value = 189534.33
$ws.Range("M70").Value = -189264.33
$ws.Range("H73").Value = 107192.37
$ws.Range("I73").Value = 189534.33
$ws.Range("K73").Value = 189534.33
$ws.Range("M73").Value = -188598.33
$ws.Range("H80").Value = 72829.375
$ws.Range("J80").Value = 3735.125
$ws.Range("L80").Value = 3735.125
$ws.Range("N80").Value = -5731.125
$ws.Range("H83").Value = 72829.375
$ws.Range("J83").Value = 3735.125
$ws.Range("L83").Value = 18675.625
$ws.Range("N83").Value = -28659.625
$ws.Range("H113").Value = 3685.75
$ws.Range("I113").Value = 2656.111
$ws.Range("J113").Value = 5539.1
$ws.Range("K113").Value = 2656.111
$ws.Range("L113").Value = 5539.1
$ws.Range("M113").Value = -486.1109999999999
$ws.Range("N113").Value = -9879.1
$ws.Range("H122").Value = 10398
$ws.Range("I122").Value = 10398
$ws.Range("K122").Value = 31194
$ws.Range("M122").Value = -28744
$ws.Range("H126").Value = 4147.7617
$ws.Range("I126").Value = 2577.7778
$ws.Range("J126").Value = 5325.25
$ws.Range("K126").Value = 7733.3334
$ws.Range("L126").Value = 15975.75
$ws.Range("M126").Value = -5263.3334
$ws.Range("N126").Value = -20915.75
$ws.Range("H132").Value = 3691.4736
$ws.Range("I132").Value = 3654.353
$ws.Range("K132").Value = 10963.059
$ws.Range("M132").Value = -8433.059000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3314.0417
$ws.Range("I7").Value = 1752.2858
$ws.Range("K7").Value = 1752.2858
$ws.Range("M7").Value = -1640.2858
$ws.Range("H40").Value = 7115.2856
$ws.Range("I40").Value = 6729.857
$ws.Range("J40").Value = 7693.4287
$ws.Range("K40").Value = 6729.857
$ws.Range("L40").Value = 7693.4287
$ws.Range("M40").Value = -6593.857
$ws.Range("N40").Value = -7965.4287
$ws.Range("H43").Value = 10000000
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H46").Value = 8610.886
$ws.Range("J46").Value = 9431.645500000001
$ws.Range("L46").Value = 9431.645500000001
$ws.Range("N46").Value = -9807.645500000001
$ws.Range("H55").Value = 6313.5
$ws.Range("J55").Value = 2502
$ws.Range("L55").Value = 2502
$ws.Range("N55").Value = -2848
$ws.Range("H61").Value = 2376.1785
$ws.Range("I61").Value = 842.4545000000001
$ws.Range("K61").Value = 842.4545000000001
$ws.Range("M61").Value = -640.4545000000001
$ws.Range("H68").Value = 5450.7
$ws.Range("I68").Value = 3898.5454
$ws.Range("K68").Value = 3898.5454
$ws.Range("M68").Value = -3149.5454
$ws.Range("H71").Value = 5450.7
$ws.Range("I71").Value = 3898.5454
$ws.Range("K71").Value = 19492.727
$ws.Range("M71").Value = -15748.727
$ws.Range("H75").Value = 9004
$ws.Range("I75").Value = 9004
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 9004
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -8068
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 9004
$ws.Range("I78").Value = 9004
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 27012
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -22332
$ws.Range("N78").ClearContents()
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996
$ws.Range("H82").Value = 3370.4546
$ws.Range("I82").Value = 1774.4667
$ws.Range("K82").Value = 1774.4667
$ws.Range("M82").Value = -1413.4667
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984
$ws.Range("H85").Value = 3370.4546
$ws.Range("I85").Value = 1774.4667
$ws.Range("K85").Value = 1774.4667
$ws.Range("M85").Value = -526.4666999999999
$ws.Range("H100").Value = 87129.84
$ws.Range("I100").Value = 159857
$ws.Range("K100").Value = 159857
$ws.Range("M100").Value = -159316
$ws.Range("H102").Value = 42599.2
$ws.Range("J102").Value = 42599.2
$ws.Range("L102").Value = 42599.2
$ws.Range("N102").Value = -49089.2
$ws.Range("H105").Value = 44038.332
$ws.Range("J105").Value = 44038.332
$ws.Range("L105").Value = 44038.332
$ws.Range("N105").Value = -51026.332
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H113").Value = 2376.1785
$ws.Range("I113").Value = 842.4545000000001
$ws.Range("K113").Value = 842.4545000000001
$ws.Range("M113").Value = 1327.5455
$ws.Range("H116").Value = 133330.08
$ws.Range("J116").Value = 133330.08
$ws.Range("L116").Value = 133330.08
$ws.Range("N116").Value = -142508.08
$ws.Range("H122").Value = 5026.5
$ws.Range("I122").Value = 4026.6843
$ws.Range("J122").Value = 7740.2856
$ws.Range("K122").Value = 12080.0529
$ws.Range("L122").Value = 23220.8568
$ws.Range("M122").Value = -9630.052899999999
$ws.Range("N122").Value = -28120.8568
$ws.Range("H126").Value = 3314.0417
$ws.Range("I126").Value = 1752.2858
$ws.Range("K126").Value = 5256.857400000001
$ws.Range("M126").Value = -2786.857400000001
$ws.Range("H132").Value = 973919.2
$ws.Range("I132").Value = 1644024
$ws.Range("J132").Value = 5990
$ws.Range("K132").Value = 4932072
$ws.Range("L132").Value = 17970
$ws.Range("M132").Value = -4929542
$ws.Range("N132").Value = -23030
$ws.Range("H136").Value = 4387.946
$ws.Range("I136").Value = 1975.9412
$ws.Range("J136").Value = 6438.15
$ws.Range("K136").Value = 5927.8236
$ws.Range("L136").Value = 19314.45
$ws.Range("M136").Value = -3377.8236
$ws.Range("N136").Value = -24414.45

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 8541.666999999999
$ws.Range("I74").Value = 4999
$ws.Range("J74").Value = 10313
$ws.Range("K74").Value = 4999
$ws.Range("L74").Value = 10313
$ws.Range("M74").Value = -4063
$ws.Range("N74").Value = -12185
$ws.Range("H75").Value = 64610
$ws.Range("I75").Value = 50000
$ws.Range("K75").Value = 50000
$ws.Range("M75").Value = -49064
$ws.Range("H77").Value = 8541.666999999999
$ws.Range("I77").Value = 4999
$ws.Range("J77").Value = 10313
$ws.Range("K77").Value = 14997
$ws.Range("L77").Value = 30939
$ws.Range("M77").Value = -10317
$ws.Range("N77").Value = -40299
$ws.Range("H78").Value = 64610
$ws.Range("I78").Value = 50000
$ws.Range("K78").Value = 150000
$ws.Range("M78").Value = -145320
$ws.Range("H81").Value = 6828.143
$ws.Range("I81").Value = 8243.1875
$ws.Range("K81").Value = 16486.375
$ws.Range("M81").Value = -15425.375
$ws.Range("H84").Value = 6828.143
$ws.Range("I84").Value = 8243.1875
$ws.Range("K84").Value = 82431.875
$ws.Range("M84").Value = -77127.875
$ws.Range("H100").Value = 721.9167
$ws.Range("I100").Value = 705.7273
$ws.Range("J100").Value = 900
$ws.Range("K100").Value = 1411.4546
$ws.Range("L100").Value = 1800
$ws.Range("M100").Value = -870.4546
$ws.Range("N100").Value = -2882
$ws.Range("H107").Value = 584.1429000000001
$ws.Range("I107").Value = 517.8
$ws.Range("K107").Value = 1553.4
$ws.Range("M107").Value = 366.6000000000001
$ws.Range("H113").Value = 438.375
$ws.Range("I113").Value = 387.5625
$ws.Range("J113").Value = 540
$ws.Range("K113").Value = 1162.6875
$ws.Range("L113").Value = 1620
$ws.Range("M113").Value = 1007.3125
$ws.Range("N113").Value = -5960
$ws.Range("H122").Value = 320182.62
$ws.Range("I122").Value = 2396.9033
$ws.Range("J122").Value = 1551602.2
$ws.Range("K122").Value = 7190.7099
$ws.Range("L122").Value = 4654806.6
$ws.Range("M122").Value = -4740.7099
$ws.Range("N122").Value = -4659706.6
$ws.Range("H132").Value = 2213.0334
$ws.Range("I132").Value = 1901.1818
$ws.Range("J132").Value = 5643.4
$ws.Range("K132").Value = 5703.5454
$ws.Range("L132").Value = 16930.2
$ws.Range("M132").Value = -3173.5454
$ws.Range("N132").Value = -21990.2
$ws.Range("H136").Value = 3118.3667
$ws.Range("I136").Value = 1700.3158
$ws.Range("J136").Value = 5567.727
$ws.Range("K136").Value = 5100.9474
$ws.Range("L136").Value = 16703.181
$ws.Range("M136").Value = -2550.9474
$ws.Range("N136").Value = -21803.181
